# Bitacora entry: "Agregue liga para administrar empleados"
#
# Appends a new log entry after the last existing paragraph
# ("Corregí detalle sobre el envío de los e-mails con mensajes vacíos."):
#   (blank)
#   Último commit: 0233bfe7b1fb80e035c92ff014ef2df68940878c
#   03/04/13
#   (blank)
#   Me quedé editando el manual del administrador general, p. 5. Estoy
#   indicando las restricciones para la nómina o nombre de usuario de
#   directores, secretarias y asistentes.
#
# The trailing "_GoBack" bookmark (Word's last-edit-position marker) must
# end up wrapping the very end of the new final paragraph instead of the
# old one.

$d = $word.ActiveDocument

# The hidden "_GoBack" bookmark currently sits at the end of the last
# paragraph. Drop it here; it gets re-created (via raw XML) at the end of
# the newly inserted content below.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Paragraph mark / run formatting shared by every entry in this log.
$pPr = '<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr></w:pPr>'
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr>'

# Blank paragraph (no run at all, matching the other spacer paragraphs in
# this document).
$blankPara = "<w:p $wNs>$pPr</w:p>"

# "Último commit: 0233bfe7b1fb80e035c92ff014ef2df68940878c"
$commitPara = "<w:p $wNs>" + $pPr +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Último </w:t></w:r>" +
    '<w:proofErr w:type="spellStart"/>' +
    "<w:r>$rPr<w:t>commit</w:t></w:r>" +
    '<w:proofErr w:type="spellEnd"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`">: </w:t></w:r>" +
    "<w:r>$rPr<w:t>0233bfe7b1fb80e035c92ff014ef2df68940878c</w:t></w:r>" +
    "</w:p>"

# "03/04/13" (renders as the top of a new page, like the rest of the log)
$datePara = "<w:p $wNs>" + $pPr +
    "<w:r>$rPr<w:lastRenderedPageBreak/><w:t>03/04/13</w:t></w:r>" +
    "</w:p>"

# Closing entry text; carries the relocated "_GoBack" bookmark at the
# very end of the document.
$entryText = 'Me quedé editando el manual del administrador general, p. 5. Estoy indicando las restricciones para la nómina o nombre de usuario de directores, secretarias y asistentes.'
$entryPara = "<w:p $wNs>" + $pPr +
    "<w:r>$rPr<w:t>$entryText</w:t></w:r>" +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    "</w:p>"

$fragment = $blankPara + $commitPara + $datePara + $blankPara + $entryPara

$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$insertionPoint.InsertXML($fragment) | Out-Null
